# Attendance workbook update:
#  - Row 2 (existing check-in record) is updated with new values
#    (student_id becomes a genuine number, name/status/elapsed-time strings change)
#  - Row 3 is appended as a brand-new attendance record
#  - Columns C, D, F get wider to fit the new content
#  - The active selection moves to F3 (last edited cell)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Row 2 : update existing record ----
$ws.Range("A2").Value = 2023130229
$ws.Range("B2").Value = "Zach"
$ws.Range("C2").Value = 45801.992698587957
$ws.Range("D2").Value = 45801.993046030089
$ws.Range("E2").Value = "PRESENT"
$ws.Range("F2").Value = "0 days 00:00:30.018904"

# ---- Row 3 : brand-new record ----
# student_id is stored as text here (matches source data), so force a
# text cell, then strip the resulting number-format override so the cell
# is left with plain default styling.
$ws.Range("A3").NumberFormat = "@"
$ws.Range("A3").Value = "2023130229"
$ws.Range("A3").ClearFormats()

$ws.Range("B3").Value = "Zach"

$ws.Range("C3").NumberFormat = "yyyy\-mm\-dd\ hh:mm:ss"
$ws.Range("C3").Value = 45802.010496557217
$ws.Range("D3").NumberFormat = "yyyy\-mm\-dd\ hh:mm:ss"
$ws.Range("D3").Value = 45802.011541173117

$ws.Range("E3").Value = "LEFT_EARLY"
$ws.Range("F3").Value = "0 days 00:01:30.254813"

# ---- column widths ----
$ws.Columns.Item(3).ColumnWidth = 22.498697916666668
$ws.Columns.Item(4).ColumnWidth = 25.830729166666668
$ws.Columns.Item(6).ColumnWidth = 19.998697916666668

# ---- selection moves to the last-entered cell ----
$ws.Range("F3").Select()
